$wb = $excel.ActiveWorkbook

# Rename "notifications" sheet to "notices"
$ws3 = $wb.Worksheets.Item("notifications")
$ws3.Name = "notices"

# Update header text in the notices sheet: "Content" -> "Contents"
$ws3.Range("B1").Value = "Contents"

# Update the body rows: "Item"/"Yes" -> "item"/"yes", "Event" -> "event", "No" -> "no"
$ws3.Range("C2").Value = "item"
$ws3.Range("D2").Value = "yes"
$ws3.Range("C3").Value = "event"
$ws3.Range("D3").Value = "no"

# Move the active selection from B4 to C4
$null = $ws3.Range("C4").Select()

# De-duplicate the identical header cell style so every header cell
# (jobs, skills, notices) resolves to the same (now single) style entry.
$ws1 = $wb.Worksheets.Item("jobs")
$ws1.Range("A1:E1").Font.Bold = $true

$ws2 = $wb.Worksheets.Item("skills")
$ws2.Range("A1:E1").Font.Bold = $true

$ws3.Range("A1:D1").Font.Bold = $true
